$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing sample counts
$ws.Range("C2").Value = 30
$ws.Range("D4").Value = 30

# Update the active selection to D5
$ws.Range("D5").Select()
